# "Add bottom of screen back/next"
# - Bumps the form_version setting (new build date)
# - Adds a new "showFooter" setting row enabling the bottom-of-screen
#   back/next navigation footer
# - Leaves the "settings" sheet selected/active (instead of
#   "table_specific_translations")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Bump form_version (settings!B3) to reflect the new build
$ws.Range("B3").Value = 20210221001

# Add the new showFooter setting in the first empty row (row 10)
$ws.Range("A10").Value = "showFooter"
$ws.Range("B10").Value = 1

# Make "settings" the active/selected sheet
$ws.Activate()
